$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''31.548.78'
$ws.Range("E2").Value = '  +5.63%  '
$ws.Range("D3").Value = '''1.708.66'
$ws.Range("E3").Value = '  +4.29%  '
$ws.Range("E4").Value = '  -0.28%  '
$ws.Range("D5").Value = '''222.02'
$ws.Range("E5").Value = '  +3.05%  '
$ws.Range("E6").Value = '  +3.10%  '
$ws.Range("E7").Value = '  -0.26%  '
$ws.Range("D8").Value = '''29.97'
$ws.Range("E8").Value = '  +4.21%  '
$ws.Range("D9").Value = '''45.33'
$ws.Range("E9").Value = '  +3.36%  '
$ws.Range("E10").Value = '  +3.75%  '
$ws.Range("E11").Value = '  +5.69%  '
$ws.Range("E12").Value = '  +1.19%  '
$ws.Range("D13").Value = '''1.953.09'
$ws.Range("E13").Value = '  +4.21%  '
$ws.Range("D14").Value = '''1.713.77'
$ws.Range("E14").Value = '  +4.48%  '
$ws.Range("D15").Value = '''10.30'
$ws.Range("E15").Value = '  +8.87%  '
$ws.Range("E16").Value = '  +3.74%  '
$ws.Range("E17").Value = '  +8.32%  '
$ws.Range("D18").Value = '''31.531.59'
$ws.Range("E18").Value = '  +5.53%  '
$ws.Range("D19").Value = '''67.23'
$ws.Range("E19").Value = '  +4.08%  '
$ws.Range("D20").Value = '''251.23'
$ws.Range("E20").Value = '  +4.53%  '
$ws.Range("E21").Value = '  +3.39%  '
$ws.Range("D22").Value = '''0.999'
$ws.Range("E22").Value = '  -0.11%  '
$ws.Range("E23").Value = '  +3.33%  '
$ws.Range("E24").Value = '  +3.14%  '
$ws.Range("E25").Value = '  -1.58%  '
$ws.Range("D26").Value = '''159.44'
$ws.Range("E26").Value = '  +1.02%  '
$ws.Range("D27").Value = '''16.06'
$ws.Range("E27").Value = '  +3.48%  '
$ws.Range("E28").Value = '  +3.11%  '
$ws.Range("D29").Value = '''6.81'
$ws.Range("E29").Value = '  +2.92%  '
$ws.Range("E30").Value = '  -0.30%  '
$ws.Range("D31").Value = '''3.78'
$ws.Range("E31").Value = '  +11.53%  '
$ws.Range("D32").Value = '''0.0505'
$ws.Range("E32").Value = '  +2.30%  '
$ws.Range("E33").Value = '  +3.89%  '
$ws.Range("E34").Value = '  +7.03%  '
$ws.Range("D35").Value = '''1.514.34'
$ws.Range("E35").Value = '  +6.26%  '
$ws.Range("E36").Value = '  +2.29%  '
$ws.Range("E37").Value = '  +2.10%  '
$ws.Range("D38").Value = '''83.64'
$ws.Range("E38").Value = '  +9.21%  '
$ws.Range("D39").Value = '''0.611'
$ws.Range("E39").Value = '  +9.14%  '
$ws.Range("E40").Value = '  +4.29%  '
$ws.Range("D41").Value = '''2.73'
$ws.Range("E41").Value = '  -0.14%  '
$ws.Range("E42").Value = '  +0.38%  '
$ws.Range("D44").Value = '''0.855'
$ws.Range("E44").Value = '  +2.49%  '
$ws.Range("E45").Value = '  +0.77%  '
$ws.Range("D46").Value = '''1.03'
$ws.Range("E46").Value = '  +2.77%  '
$ws.Range("E47").Value = '  -0.19%  '
$ws.Range("D48").Value = '''52.47'
$ws.Range("E48").Value = '  +7.40%  '
$ws.Range("D49").Value = '''5.58'
$ws.Range("E49").Value = '  +3.31%  '
$ws.Range("D50").Value = '''1.841.62'
$ws.Range("E50").Value = '  +3.39%  '
$ws.Range("E51").Value = '  +9.62%  '
